$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-19 Wednesday" "2025-03-20 Thursday"

Replace-Text "476÷9=" "908÷4="
Replace-Text "471÷7=" "625÷5="
Replace-Text "590÷7=" "689÷2="
Replace-Text "986÷8=" "434÷2="
Replace-Text "651÷4=" "391÷7="

Replace-Text "510÷7=" "481÷2="
Replace-Text "353÷7=" "824÷7="
Replace-Text "251÷7=" "150÷3="
Replace-Text "958÷3=" "831÷6="
Replace-Text "147÷8=" "562÷6="

Replace-Text "355÷3=" "874÷4="
Replace-Text "869÷3=" "591÷4="
Replace-Text "232÷6=" "859÷7="
Replace-Text "940÷5=" "980÷6="
Replace-Text "812÷3=" "444÷9="

Replace-Text "844÷6=" "503÷7="
Replace-Text "598÷5=" "447÷7="
Replace-Text "634÷9=" "651÷5="
Replace-Text "900÷4=" "806÷2="
Replace-Text "542÷7=" "490÷2="

Replace-Text "529÷9=" "383÷7="
Replace-Text "956÷3=" "412÷4="
Replace-Text "849÷4=" "909÷6="
Replace-Text "968÷4=" "138÷4="
Replace-Text "894÷2=" "964÷9="
